$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of metric data (row 87)
$ws.Cells.Item(87, 1).Value = "2025-04-29 14:32:28"
$ws.Cells.Item(87, 2).Value = 244
